$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value corrections on existing rows ---
$ws.Range("B3").Value = 11
$ws.Range("E5").Value = 4
$ws.Range("Q5").Value = 10
$ws.Range("H26").Value = 1
$ws.Range("T30").Value = 1
$ws.Range("H35").Value = 1
$ws.Range("N49").Value = 3
$ws.Range("E60").Value = 1
$ws.Range("E62").Value = 2

# --- Insert a new neighborhood row (O'Hare) before the current row 75 (Oakland) ---
# This shifts Oakland..Wrigleyville down by one row, to rows 76..91.
$ws.Rows("75:75").Insert()

# Match the label style used by the other neighborhood-name cells in column A
$ws.Range("A74").Copy() | Out-Null
$ws.Range("A75").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A75").Value = "O'Hare"
$ws.Range("H75").Value = 1

# Add one more newly-observed data point for Boystown (now row 78 after the shift)
$ws.Range("H78").Value = 1

# --- Add the brand new trailing row for Wrigleyville's data is already carried by the
#     row-insert shift above (old row 90 -> new row 91); nothing else to do there. ---

# --- Header / label text updates ---
$ws.Range("B1").Value = "March 2022 (through March 26)"
$ws.Name = "Through 2022-03-26"
